$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Adapt spreadsheet inputs for 6V supercap (user input cells only;
# all dependent formulas recalculate automatically).
$ws.Range("I8").Value = 1.3
$ws.Range("C9").Value = 5.5
$ws.Range("I9").Value = 1.5
$ws.Range("O9").Value = 3.5
$ws.Range("C19").Value = 4.32
$ws.Range("I19").Value = 10
$ws.Range("O19").Value = 4.53
$ws.Range("C20").Value = 8.87
$ws.Range("I20").Value = 0.787
$ws.Range("I21").Value = 1.74
$ws.Range("O21").Value = 8.66

# Reflect the cell selection that was active when the workbook was saved.
$ws.Range("I21").Select()
